$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list): update "想去人数" (column F) for rows 3-10
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 263
$ws1.Range("F4").Value = 147
$ws1.Range("F5").Value = 1776
$ws1.Range("F6").Value = 1525
$ws1.Range("F7").Value = 277
$ws1.Range("F8").Value = 65
$ws1.Range("F9").Value = 495
$ws1.Range("F10").Value = 124

# Sheet "全部类型" (all types list): update "想去人数" (column F) for rows 3-7 and 9-11
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 263
$ws4.Range("F4").Value = 147
$ws4.Range("F5").Value = 1776
$ws4.Range("F6").Value = 1525
$ws4.Range("F7").Value = 277
$ws4.Range("F9").Value = 65
$ws4.Range("F10").Value = 495
$ws4.Range("F11").Value = 124

$wb.Save()
